$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SJU", "Puerto Rico"),
    @("BGI", "Barbados"),
    @("SDQ", "Dominican Republic"),
    @("CUN", "Cancun"),
    @("SXM", "Saint Maarten"),
    @("AUA", "Aruba"),
    @("UVF", "Saint Lucia"),
    @("NAS", "Bahamas"),
    @("CUR", "Curacao"),
    @("PLS", "Providenciales"),
    @("MBJ", "Jamaica"),
    @("LIR", "Costa Rica"),
    @("BDA", "Bermuda"),
    @("GCM", "Cayman Islands"),
    @("GND", "Grenada"),
    @("POP", "Dominican Republic"),
    @("PUJ", "Dominican Republic"),
    @("AZS", "Dominican Republic"),
    @("LRM", "Dominican Republic")
)

$startRow = 188
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}

$ws.Range("D205").Select()
